$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("C2").Value = 7.1
$ws.Range("B5").Value = 0.95

# Update the active selection from C5 to B5
$ws.Range("B5").Select()
